# Scheduled runner update: refresh market price / profit columns (H-N)
# across the FFXIV Leve Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 950.375
$ws.Range("I2").Value = 871.8570999999999
$ws.Range("K2").Value = 871.8570999999999
$ws.Range("M2").Value = -758.8570999999999

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 335.8
$ws.Range("I41").Value = 335.8
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 335.8
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 104.2
$ws.Range("N41").ClearContents()

# ALC row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 221.16667
$ws.Range("I53").Value = 224.5
$ws.Range("K53").Value = 224.5
$ws.Range("M53").Value = 412.5

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5964.6665
$ws.Range("I76").Value = 5964.6665
$ws.Range("K76").Value = 5964.6665
$ws.Range("M76").Value = -5649.6665

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5964.6665
$ws.Range("I79").Value = 5964.6665
$ws.Range("K79").Value = 5964.6665
$ws.Range("M79").Value = -4872.6665

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7252
$ws.Range("I86").Value = 7252
$ws.Range("K86").Value = 7252
$ws.Range("M86").Value = -6129

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7252
$ws.Range("I89").Value = 7252
$ws.Range("K89").Value = 36260
$ws.Range("M89").Value = -30644

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2002
$ws.Range("I92").Value = 2002
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2002
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -754
$ws.Range("N92").ClearContents()

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 8241.053
$ws.Range("I98").Value = 2440
$ws.Range("K98").Value = 2440
$ws.Range("M98").Value = -942

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = 254

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 8241.053
$ws.Range("I122").Value = 2440
$ws.Range("K122").Value = 7320
$ws.Range("M122").Value = -4870

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 975
$ws.Range("I141").Value = 975
$ws.Range("K141").Value = 2925
$ws.Range("M141").Value = 2255

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 410.1111
$ws.Range("I32").Value = 410.1111
$ws.Range("K32").Value = 410.1111
$ws.Range("M32").Value = -123.1111

# BSM row 37
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

# BSM row 125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

# BSM row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 40750
$ws.Range("J130").Value = 40750
$ws.Range("L130").Value = 40750
$ws.Range("N130").Value = -50790

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 432.5
$ws.Range("I16").Value = 432.5
$ws.Range("K16").Value = 432.5
$ws.Range("M16").Value = -145.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5542.1
$ws.Range("I31").Value = 4781.1333
$ws.Range("K31").Value = 4781.1333
$ws.Range("M31").Value = -4486.1333

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5542.1
$ws.Range("I34").Value = 4781.1333
$ws.Range("K34").Value = 4781.1333
$ws.Range("M34").Value = -4579.1333

# CRP row 103
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

# CRP row 104
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -45242

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 432.5
$ws.Range("I113").Value = 432.5
$ws.Range("K113").Value = 432.5
$ws.Range("M113").Value = 1737.5

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.916668
$ws.Range("I12").Value = 35.5
$ws.Range("J12").Value = 30.333334
$ws.Range("K12").Value = 106.5
$ws.Range("L12").Value = 91.00000199999999
$ws.Range("M12").Value = 66.5
$ws.Range("N12").Value = -437.000002

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 126
$ws.Range("J38").Value = 172
$ws.Range("L38").Value = 516
$ws.Range("N38").Value = -1210

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14400
$ws.Range("I131").Value = 38618
$ws.Range("J131").Value = 2291
$ws.Range("K131").Value = 115854
$ws.Range("L131").Value = 6873
$ws.Range("M131").Value = -110814
$ws.Range("N131").Value = -16953

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3002
$ws.Range("N80").ClearContents()

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -15008
$ws.Range("N83").ClearContents()

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1492.5714
$ws.Range("I107").Value = 1492.5714
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1492.5714
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 427.4286
$ws.Range("N107").ClearContents()

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9999
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 3000
$ws.Range("I48").Value = 3000
$ws.Range("K48").Value = 3000
$ws.Range("M48").Value = -2339

# LTW row 121
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# LTW row 130
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 100427
$ws.Range("J130").Value = 100427
$ws.Range("L130").Value = 100427
$ws.Range("N130").Value = -110467

# WVR row 131
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 93749.75
$ws.Range("I131").Value = 100000
$ws.Range("J131").Value = 91666.336
$ws.Range("K131").Value = 100000
$ws.Range("L131").Value = 91666.336
$ws.Range("M131").Value = -94960
$ws.Range("N131").Value = -101746.336

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5559.9565
$ws.Range("I136").Value = 3615.4167
$ws.Range("J136").Value = 7681.273
$ws.Range("K136").Value = 10846.2501
$ws.Range("L136").Value = 23043.819
$ws.Range("M136").Value = -8296.250100000001
$ws.Range("N136").Value = -28143.819

# WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1900
$ws.Range("J14").Value = 1900
$ws.Range("L14").Value = 1900
$ws.Range("N14").Value = -2236

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9093.5625
$ws.Range("I132").Value = 5866.3335
$ws.Range("J132").Value = 11029.9
$ws.Range("K132").Value = 17599.0005
$ws.Range("L132").Value = 33089.7
$ws.Range("M132").Value = -15069.0005
$ws.Range("N132").Value = -38149.7

